$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITR input data")
$ws.Activate()
Write-Host "SplitColumn:" $excel.ActiveWindow.SplitColumn
Write-Host "SplitRow:" $excel.ActiveWindow.SplitRow
Write-Host "Split:" $excel.ActiveWindow.Split
